$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns keep their text (string) type when we
# write numeric-looking values into them (Excel would otherwise auto-convert
# strings like "1.001" or "0.9999" into numbers).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.291.02"
$ws.Range("E2").Value = "  +0.65%  "
$ws.Range("D3").Value = "1.903.24"
$ws.Range("E3").Value = "  +0.44%  "
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").Value = "306.61"
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("D7").Value = "0.5425"
$ws.Range("E7").Value = "  +4.11%  "
$ws.Range("D8").Value = "0.3811"
$ws.Range("E8").Value = "  +1.25%  "
$ws.Range("D9").Value = "0.07309"
$ws.Range("E9").Value = "  +0.64%  "
$ws.Range("D10").Value = "22.09"
$ws.Range("E10").Value = "  +4.39%  "
$ws.Range("D11").Value = "0.9025"
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("D12").Value = "0.08189"
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("D13").Value = "95.55"
$ws.Range("E13").Value = "  -0.69%  "
$ws.Range("D14").Value = "5.352"
$ws.Range("E14").Value = "  +0.92%  "
$ws.Range("D15").Value = "0.9994"
$ws.Range("E15").Value = "  -0.24%  "
$ws.Range("D16").Value = "14.89"
$ws.Range("E16").Value = "  +2.09%  "
$ws.Range("D17").Value = "0.000008659"
$ws.Range("E17").Value = "  +0.70%  "
$ws.Range("D18").Value = "1.378.57"
$ws.Range("E18").Value = "  -27.63%  "
$ws.Range("D19").Value = "0.9985"
$ws.Range("E19").Value = "  -0.48%  "
$ws.Range("D20").Value = "27.298.50"
$ws.Range("E20").Value = "  +0.51%  "
$ws.Range("D21").Value = "5.057"
$ws.Range("E21").Value = "  -0.55%  "
$ws.Range("D22").Value = "10.82"
$ws.Range("E22").Value = "  +1.15%  "
$ws.Range("D23").Value = "6.515"
$ws.Range("E23").Value = "  +1.57%  "
$ws.Range("D24").Value = "148.55"
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("D25").Value = "2.312"
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("D26").Value = "18.35"
$ws.Range("E26").Value = "  +0.87%  "
$ws.Range("D27").Value = "1.756"
$ws.Range("E27").Value = "  +1.35%  "
$ws.Range("E28").Value = "  +1.23%  "
$ws.Range("D29").Value = "4.850"
$ws.Range("E29").Value = "  +1.19%  "
$ws.Range("E30").Value = "  -3.88%  "
$ws.Range("D31").Value = "0.09201"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D32").Value = "0.8264"
$ws.Range("E32").Value = "  +4.23%  "
$ws.Range("D33").Value = "0.05077"
$ws.Range("E33").Value = "  +1.00%  "
$ws.Range("D34").Value = "1.225"
$ws.Range("E34").Value = "  +0.70%  "
$ws.Range("D35").Value = "3.011"
$ws.Range("E35").Value = "  +1.16%  "
$ws.Range("D36").Value = "3.318"
$ws.Range("E36").Value = "  -3.28%  "
$ws.Range("E37").Value = "  +3.03%  "
$ws.Range("D38").Value = "0.5995"
$ws.Range("E38").Value = "  +4.77%  "
$ws.Range("D39").Value = "0.01997"
$ws.Range("E39").Value = "  +0.11%  "
$ws.Range("E40").Value = "  +0.34%  "
$ws.Range("D41").Value = "9.285"
$ws.Range("E41").Value = "  +2.87%  "
$ws.Range("D42").Value = "6.671"
$ws.Range("E42").Value = "  +1.63%  "
$ws.Range("D43").Value = "115.97"
$ws.Range("E43").Value = "  -0.39%  "
$ws.Range("D44").Value = "0.5153"
$ws.Range("E44").Value = "  +5.99%  "
$ws.Range("D45").Value = "0.1531"
$ws.Range("E45").Value = "  +1.11%  "
$ws.Range("D46").Value = "10.22"
$ws.Range("E46").Value = "  +1.84%  "
$ws.Range("D47").Value = "1.001"
$ws.Range("E47").Value = "  -0.15%  "
$ws.Range("D48").Value = "1.639"
$ws.Range("E48").Value = "  +1.16%  "
$ws.Range("D49").Value = "38.12"
$ws.Range("E49").Value = "  -0.30%  "
$ws.Range("D50").Value = "0.06097"
$ws.Range("E50").Value = "  +2.89%  "
$ws.Range("D51").Value = "63.42"
$ws.Range("E51").Value = "  -0.26%  "
